$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.128865
$ws.Range("H2").Value = 69.386595
$ws.Range("I2").Value = 0.7917836846260858
$ws.Range("J2").Value = 0.7917836846260858
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 2602.30705426008
$ws.Range("R2").Value = 23420.76348834072
$ws.Range("S2").Value = 0.2593339377524517
$ws.Range("T2").Value = 0.2593339377524517

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.128865
$ws.Range("H3").Value = 69.386595
$ws.Range("I3").Value = 0.7917836846260858
$ws.Range("J3").Value = 0.7917836846260858
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 2458.93293166109
$ws.Range("R3").Value = 22130.39638494981
$ws.Range("S3").Value = 0.2450459329128495
$ws.Range("T3").Value = 0.2450459329128496

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.128865
$ws.Range("H4").Value = 69.386595
$ws.Range("I4").Value = 0.7917836846260858
$ws.Range("J4").Value = 0.7917836846260858
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 2883.976462831197
$ws.Range("R4").Value = 25955.78816548077
$ws.Range("S4").Value = 0.2874038139607846
$ws.Range("T4").Value = 0.2874038139607846

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.158105666666667
$ws.Range("H5").Value = 3.474317
$ws.Range("I5").Value = 0.03964609469334889
$ws.Range("J5").Value = 0.03964609469334889
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 130.302396851088
$ws.Range("R5").Value = 1172.721571659792
$ws.Range("S5").Value = 0.01298533684511086
$ws.Range("T5").Value = 0.01298533684511086

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.158105666666667
$ws.Range("H6").Value = 3.474317
$ws.Range("I6").Value = 0.03964609469334889
$ws.Range("J6").Value = 0.03964609469334889
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 123.1233855232407
$ws.Range("R6").Value = 1108.110469709166
$ws.Range("S6").Value = 0.0122699096345623
$ws.Range("T6").Value = 0.01226990963456231

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.158105666666667
$ws.Range("H7").Value = 3.474317
$ws.Range("I7").Value = 0.03964609469334889
$ws.Range("J7").Value = 0.03964609469334889
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 144.4061126275802
$ws.Range("R7").Value = 1299.655013648222
$ws.Range("S7").Value = 0.01439084821367573
$ws.Range("T7").Value = 0.01439084821367573

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.924119999999999
$ws.Range("H8").Value = 14.77236
$ws.Range("I8").Value = 0.1685702206805652
$ws.Range("J8").Value = 0.1685702206805652
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 554.0294438150399
$ws.Range("R8").Value = 4986.264994335359
$ws.Range("S8").Value = 0.05521202313929379
$ws.Range("T8").Value = 0.05521202313929379

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.924119999999999
$ws.Range("H9").Value = 14.77236
$ws.Range("I9").Value = 0.1685702206805652
$ws.Range("J9").Value = 0.1685702206805652
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 523.50518831992
$ws.Range("R9").Value = 4711.54669487928
$ws.Range("S9").Value = 0.05217011639675447
$ws.Range("T9").Value = 0.05217011639675448

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.924119999999999
$ws.Range("H10").Value = 14.77236
$ws.Range("I10").Value = 0.1685702206805652
$ws.Range("J10").Value = 0.1685702206805652
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 613.9966738599733
$ws.Range("R10").Value = 5525.97006473976
$ws.Range("S10").Value = 0.06118808114451697
$ws.Range("T10").Value = 0.06118808114451698
